$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain text values (dotted thousands-separator
# notation such as "25.694.87", trailing-zero decimals such as "6.800", etc.).
# Pin every Price cell we touch to the Text number format before writing it so
# Excel stores the exact refreshed digits instead of re-parsing them as a number
# (and dropping e.g. a trailing zero). Cells we do not touch are left untouched.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.694.87"
$ws.Range("E2").Value = "  -3.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.631.11"
$ws.Range("E3").Value = "  -2.35%  "
$ws.Range("E4").Value = "  +0.97%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.88"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5043"
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.014"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2572"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06382"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.45"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07737"
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.673.43"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.241"
$ws.Range("E13").Value = "  -2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.853.40"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5446"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7880"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.21"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.739.39"
$ws.Range("E18").Value = "  -3.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.013"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "204.01"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.334"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.961"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.926"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.015"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.966"
$ws.Range("E25").Value = "  +14.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.75"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1152"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.71"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.800"
$ws.Range("E29").Value = "  -2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.242"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05008"
$ws.Range("E31").Value = "  -3.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.257"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.197"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.533"
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.349"
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.649"
$ws.Range("E36").Value = "  -3.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8962"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5650"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.120.80"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01558"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.590"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.017"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.640"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8106"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.42"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.766.94"
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4538"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.013"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.71"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05047"
$ws.Range("E51").Value = "  -1.76%  "
